# Remove VLOOKUP from excel
# 1) Delete the helper lookup sheets (MileageBand, PriceBands)
# 2) On VehicleData, drop the "Year Filter" column and the columns that
#    depended on VLOOKUPs into the deleted sheets (Price Band, Mileage Band)
# 3) Replace the remaining formulas with their static computed values and
#    relabel the shifted header columns

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("VehicleData")

# --- Drop helper columns that are no longer needed ------------------------
# W: Price Band (VLOOKUP into PriceBands)   -> remove
# Y: Mileage Band (VLOOKUP into MileageBand) -> remove
# AD: Year Filter                            -> remove
# Delete right-to-left so earlier column letters stay valid.
$ws.Columns("AD").Delete()
$ws.Columns("Y").Delete()
$ws.Columns("W").Delete()

# Now delete the now-unused lookup sheets
$wb.Worksheets("MileageBand").Delete()
$wb.Worksheets("PriceBands").Delete()

# --- Update header row ------------------------------------------------------
$ws.Range("U1").Value = "Transmission Mapped"
$ws.Range("V1").Value = "Value Rounded to Nearest 5000"
$ws.Range("W1").Value = "Mileage Rounded to Nearest 50,000"
$ws.Range("X1").Value = "Engine Size Rounded"
$ws.Range("Y1").Value = "Price Filter"
$ws.Range("Z1").Value = "Mileage Filter"
$ws.Range("AA1").Value = "Engine Size Filter"
$ws.Range("AB1").Value = "MPG Filter"
$ws.Range("AC1").Value = "Master Filter"

# --- Replace formulas with static values (row 2) ---------------------------
$ws.Range("U2").Value = "Automatic"
$ws.Range("V2").Value = 10000
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1.6
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 1
$ws.Range("AA2").Value = 1
$ws.Range("AB2").Value = 1
$ws.Range("AC2").Value = 1

# --- Replace formulas with static values (row 3) ---------------------------
$ws.Range("U3").Value = "Manual"
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 1
$ws.Range("Y3").Value = 1
$ws.Range("Z3").Value = 1
$ws.Range("AA3").Value = 1
$ws.Range("AB3").Value = 1
$ws.Range("AC3").Value = 1
